# RQSD RPS Qualifying Source Definitions.xlsx - "Initial update with CPL's work-to-date"
#
# Summary of content-level changes applied:
#  1. "About" sheet: remove the note row that referenced the
#     "We are following the EUA model ..." text (the row + the blank spacer
#     row right after it), which shifts all subsequent notes rows up and
#     lets the now-unused shared string drop out of the table on save.
#  2. "RQSD-BRQSD" sheet: the special highlight formatting on the
#     nuclear / onshore-wind / municipal-solid-waste rows is removed, and
#     the "Qualifies for RPS" flag for municipal solid waste flips 1 -> 0.
#  3. "RQSD-RQSD" sheet: the highlight formatting that covered the whole
#     data block (rows 4-17) is removed, and the municipal-solid-waste
#     flag flips 1 -> 0 there as well.
#  4. Selection bookmarks on the two RQSD sheets move to B2 (matching the
#     state the workbook was left in), while the About sheet stays the
#     active/selected tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "About" sheet — drop the EUA-model note row (and the blank row that
#    trails it), shifting every later row up by two.
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Rows("9:10").Delete()

# ---------------------------------------------------------------------
# 2. "RQSD-BRQSD" sheet — clear the special cell formatting that had been
#    applied to the nuclear, onshore wind, and municipal solid waste rows,
#    and flip the municipal solid waste RPS flag to 0.
# ---------------------------------------------------------------------
$wsBrqsd = $wb.Worksheets.Item("RQSD-BRQSD")
$wsBrqsd.Range("A4:B4").ClearFormats()
$wsBrqsd.Range("A10:B10").ClearFormats()
$wsBrqsd.Range("A17:B17").ClearFormats()
$wsBrqsd.Range("B17").Value = 0

# ---------------------------------------------------------------------
# 3. "RQSD-RQSD" sheet — clear the special formatting across the whole
#    data block, and flip the municipal solid waste RPS flag to 0.
# ---------------------------------------------------------------------
$wsRqsd = $wb.Worksheets.Item("RQSD-RQSD")
$wsRqsd.Range("A4:B17").ClearFormats()
$wsRqsd.Range("B17").Value = 0

# ---------------------------------------------------------------------
# 4. Leave the cell-pointer bookmark on each RQSD sheet at B2, then
#    return focus to the About sheet so it remains the selected tab.
# ---------------------------------------------------------------------
$wsBrqsd.Activate()
$wsBrqsd.Range("B2").Select() | Out-Null

$wsRqsd.Activate()
$wsRqsd.Range("B2").Select() | Out-Null

$wsAbout.Activate()
